# Reorder the keys inside the dictionary-like strings in column A of the
# "Subgroups" worksheet so that they match the new canonical ordering
# produced by the author's updated algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$updates = @{
    10 = "{'Hobby': '1', 'Student': '1'}"
    11 = "{'Hobby': '1', 'SexualOrientation': '1'}"
    13 = "{'Dependents': '2', 'Hobby': '1'}"
    17 = "{'HDI': '1', 'Student': '1'}"
    18 = "{'UndergradMajor': '2', 'SexualOrientation': '1'}"
    21 = "{'Dependents': '2', 'SexualOrientation': '1'}"
    22 = "{'HDI': '1', 'SexualOrientation': '1'}"
    24 = "{'Hobby': '1', 'Student': '1', 'SexualOrientation': '1'}"
    25 = "{'RaceEthnicity': '1', 'Hobby': '1', 'SexualOrientation': '1'}"
    26 = "{'Hobby': '1', 'SexualOrientation': '1', 'HDI': '1'}"
    28 = "{'HDI': '1', 'Student': '1', 'SexualOrientation': '1'}"
    29 = "{'RaceEthnicity': '1', 'HDI': '1', 'Student': '1'}"
    30 = "{'RaceEthnicity': '1', 'HDI': '1', 'SexualOrientation': '1'}"
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
